$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted as row 25 ("Fruta / hortaliza, semanal").
# Every existing row from 25 downward (previously 25..106) shifts down by one
# (to 26..107); the sheet dimension grows from R106 to R107 automatically.
$ws.Rows.Item(25).Insert()

$ws.Cells.Item(25, 1).Value = 11
$ws.Cells.Item(25, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(25, 3).Value = "Bíobío"
$ws.Cells.Item(25, 4).Value = 45030
$ws.Cells.Item(25, 5).Value = 8
$ws.Cells.Item(25, 6).Value = 100112037
$ws.Cells.Item(25, 7).Value = "Cebollín"
$ws.Cells.Item(25, 8).Value = "Sin especificar"
$ws.Cells.Item(25, 9).Value = "Primera"
$ws.Cells.Item(25, 10).Value = 130
$ws.Cells.Item(25, 11).Value = 5000
$ws.Cells.Item(25, 12).Value = 6000
$ws.Cells.Item(25, 13).Value = 5385
$ws.Cells.Item(25, 14).Value = "`$/paquete 36 unidades"
$ws.Cells.Item(25, 15).Value = "Región Metropolitana"
$ws.Cells.Item(25, 16).Value = 150
$ws.Cells.Item(25, 17).Value = 36
$ws.Cells.Item(25, 18).Value = "Hortaliza"
